# Build the "schedule_template" worksheet: a small header block (Date /
# Day-of-week) plus a table header row and an empty 12-row schedule grid.
#
# Excel/COM constants used below:
#   xlThin = 2            xlMedium = -4138
#   xlEdgeLeft = 7  xlEdgeTop = 8  xlEdgeBottom = 9  xlEdgeRight = 10
#   xlLeft = -4131  xlCenter = -4108

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell values, written in reading order so the shared-string table
#    comes out in the same order as the source workbook.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Дата"
$ws.Range("B1").Value = "День недели"
$ws.Range("B2").Value = "Понедельник"
$ws.Range("A3").Value = "Цвет"
$ws.Range("B3").Value = "Время начала"
$ws.Range("C3").Value = "Кабинет"
$ws.Range("D3").Value = "Сод-ние"
$ws.Range("E3").Value = "!"
$ws.Range("F3").Value = "Карта"
$ws.Range("G3").Value = "Врач"
$ws.Range("H3").Value = "Пациент"

# A2 holds the actual date serial (2023-02-27, a Monday) - numeric, not text.
$ws.Range("A2").Value = 44984

# ---------------------------------------------------------------------
# 2. Borders + center/center alignment for the header cells. Borders are
#    set per single cell (not per range) so that left/right edges land on
#    every cell instead of only the outer edge of a multi-cell range.
#    NumberFormat is applied last (see Set-CellBox) - applying it before
#    alignment on a freshly bordered cell makes the host re-derive the
#    format string and lose the built-in numFmtId 14 mapping.
# ---------------------------------------------------------------------
function Set-CellBox($cell, $left, $right, $top, $bottom, $dateFmt) {
    if ($left -ne 0)   { $cell.Borders.Item(7).Weight = $left }
    if ($right -ne 0)  { $cell.Borders.Item(10).Weight = $right }
    if ($top -ne 0)    { $cell.Borders.Item(8).Weight = $top }
    if ($bottom -ne 0) { $cell.Borders.Item(9).Weight = $bottom }
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
    if ($dateFmt) { $cell.NumberFormat = "mm-dd-yy" }
}

$xlThin = 2
$xlMedium = -4138

# Row 3 table header: A3 + H3 are the thick-boxed ends, B3:G3 share a thin
# divider between every column.
Set-CellBox $ws.Range("A3") $xlMedium $xlThin  $xlMedium $xlMedium $false
foreach ($col in 2..7) {
    Set-CellBox $ws.Cells.Item(3, $col) $xlThin $xlThin $xlMedium $xlMedium $false
}
Set-CellBox $ws.Range("H3") $xlThin $xlMedium $xlMedium $xlMedium $false

# Row 1 / Row 2: small "Date" / "Day of week" box above the table.
Set-CellBox $ws.Range("A1") $xlMedium $xlThin  $xlMedium $xlMedium $false
Set-CellBox $ws.Range("B1") 0         $xlMedium $xlMedium $xlMedium $false
Set-CellBox $ws.Range("A2") $xlMedium $xlThin  0          0         $true
Set-CellBox $ws.Range("B2") 0         $xlMedium 0          0        $false

# ---------------------------------------------------------------------
# 3. Alignment (then number format) for the (empty) C1:G1 run - mirrors
#    the date-format cell above but without borders.
# ---------------------------------------------------------------------
$ws.Range("C1:G1").HorizontalAlignment = -4131
$ws.Range("C1:G1").VerticalAlignment = -4108
$ws.Range("C1:G1").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------
# 4. Empty schedule body (rows 4-15, columns A-H): left/center aligned,
#    no border.
# ---------------------------------------------------------------------
$ws.Range("A4:H15").HorizontalAlignment = -4131
$ws.Range("A4:H15").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 5. Column widths (character units), approximating the source sheet.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 13.7109375
$ws.Columns.Item(3).ColumnWidth = 8
$ws.Columns.Item(4).ColumnWidth = 8.7109375
$ws.Columns.Item(5).ColumnWidth = 5.7109375
$ws.Columns.Item(6).ColumnWidth = 13.7109375
$ws.Columns.Item(7).ColumnWidth = 30.7109375
$ws.Columns.Item(8).ColumnWidth = 30.7109375

# ---------------------------------------------------------------------
# 6. Selection, matching the saved view state.
# ---------------------------------------------------------------------
[void]$ws.Range("A2").Select()

Write-Host "schedule_template layout applied"
